$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("D5", "D8", "D10", "D11", "D16", "D18", "D20", "D23", "D25", "D27", "D32", "D36", "D38", "D39", "D40", "D44", "D45", "D48", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.618.97"
$ws.Range("E2").Value = "  -2.57%  "
$ws.Range("D3").Value = "1.665.84"
$ws.Range("E3").Value = "  -4.08%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "215.58"
$ws.Range("E5").Value = "  -2.11%  "
$ws.Range("E6").Value = "  -2.60%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "24.10"
$ws.Range("E8").Value = "  -0.78%  "
$ws.Range("E9").Value = "  -1.12%  "
$ws.Range("D10").Value = "0.0620"
$ws.Range("E10").Value = "  -2.58%  "
$ws.Range("D11").Value = "0.0878"
$ws.Range("E11").Value = "  -1.96%  "
$ws.Range("D12").Value = "1.902.53"
$ws.Range("E12").Value = "  -4.02%  "
$ws.Range("D13").Value = "1.642.01"
$ws.Range("E13").Value = "  -5.40%  "
$ws.Range("E14").Value = "  -3.40%  "
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("D16").Value = "66.47"
$ws.Range("E16").Value = "  -1.99%  "
$ws.Range("D17").Value = "27.642.52"
$ws.Range("E17").Value = "  -2.48%  "
$ws.Range("D18").Value = "241.86"
$ws.Range("E18").Value = "  -0.66%  "
$ws.Range("E19").Value = "  -3.58%  "
$ws.Range("D20").Value = "7.69"
$ws.Range("E20").Value = "  -4.13%  "
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("E22").Value = "  -3.35%  "
$ws.Range("D23").Value = "9.37"
$ws.Range("E23").Value = "  -3.51%  "
$ws.Range("E24").Value = "  -3.46%  "
$ws.Range("D25").Value = "147.35"
$ws.Range("E25").Value = "  -1.64%  "
$ws.Range("E26").Value = "  -4.01%  "
$ws.Range("D27").Value = "16.40"
$ws.Range("E27").Value = "  -2.22%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("E29").Value = "  -2.58%  "
$ws.Range("E30").Value = "  +1.60%  "
$ws.Range("E31").Value = "  -1.99%  "
$ws.Range("D32").Value = "3.35"
$ws.Range("E32").Value = "  -2.73%  "
$ws.Range("D33").Value = "1.463.39"
$ws.Range("E33").Value = "  -2.81%  "
$ws.Range("E34").Value = "  -4.26%  "
$ws.Range("E35").Value = "  -4.89%  "
$ws.Range("D36").Value = "0.930"
$ws.Range("E36").Value = "  -4.14%  "
$ws.Range("E37").Value = "  -1.31%  "
$ws.Range("D38").Value = "0.579"
$ws.Range("E38").Value = "  -4.77%  "
$ws.Range("D39").Value = "0.0172"
$ws.Range("E39").Value = "  -2.11%  "
$ws.Range("D40").Value = "70.00"
$ws.Range("E40").Value = "  -1.21%  "
$ws.Range("E41").Value = "  -4.31%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("E43").Value = "  -3.58%  "
$ws.Range("D44").Value = "5.41"
$ws.Range("E44").Value = "  -5.63%  "
$ws.Range("D45").Value = "0.795"
$ws.Range("E45").Value = "  -0.94%  "
$ws.Range("D46").Value = "1.810.26"
$ws.Range("E46").Value = "  -3.94%  "
$ws.Range("E47").Value = "  +0.87%  "
$ws.Range("D48").Value = "88.94"
$ws.Range("E48").Value = "  -2.54%  "
$ws.Range("D49").Value = "0.0₆0108"
$ws.Range("E49").Value = "  -5.02%  "
$ws.Range("E50").Value = "  -2.32%  "
$ws.Range("D51").Value = "7.93"
$ws.Range("E51").Value = "  -3.60%  "

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
